# Applies the "mit_ettem" (what-did-I-eat) tracker update:
#  - a few portion/quantity corrections
#  - a new "sertés" (pork) meal logged on row 4 (columns I-L)
#  - the dates on rows 9-13 shifted forward so each day is unique
#  - a new "saláta" (salad) portion logged on row 13 (columns I-J)
#  - cosmetic: active selection moved to G19, column I widened

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: breakfast portion size corrected 200 -> 100 ---
$ws.Range("C2").Value = 100

# --- Row 3: lunch portion size corrected 30 -> 300 ---
$ws.Range("E3").Value = 300

# --- Row 4: new 4th meal added (sertés 300 / rizs 400) ---
$ws.Range("I4").Value = 300
$ws.Range("J4").Value = "sertés"
$ws.Range("K4").Value = 400
$ws.Range("L4").Value = "rizs"

# --- Rows 9-13: dates de-duplicated / shifted forward one day each ---
$ws.Range("A9").Value  = 45918
$ws.Range("A10").Value = 45919
$ws.Range("A11").Value = 45920
$ws.Range("A12").Value = 45921
$ws.Range("A13").Value = 45922

# --- Row 13: new 4th meal added (saláta 300) ---
$ws.Range("I13").Value = 300
$ws.Range("J13").Value = "saláta"

# --- Cosmetic: widen column I slightly to fit the new entries ---
$ws.Columns.Item(9).ColumnWidth = 12.436197916666666

# --- Cosmetic: move the active selection / window focus to G19 ---
$ws.Range("G19").Select()

# --- Cosmetic: shrink/reposition the workbook window (best effort;
#     some hosts treat window chrome as session-only state) ---
$win = $wb.Windows.Item(1)
$win.Left   = -110
$win.Top    = -110
$win.Width  = 19420
$win.Height = 10300
